$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.708.64"
$ws.Range("E2").Value = "  -0.54%  "

$ws.Range("D3").Value = "2.520.02"
$ws.Range("E3").Value = "  -1.29%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "315.57"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.53%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "95.29"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.02%  "

$ws.Range("E7").Value = "  +0.30%  "

$ws.Range("E8").Value = "  +0.02%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.538"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.56%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.20"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.24%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0809"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.41%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.55"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.62%  "

$ws.Range("E13").Value = "  -2.06%  "

$ws.Range("D14").Value = "2.908.24"
$ws.Range("E14").Value = "  -1.22%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.60"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +4.00%  "

$ws.Range("D16").Value = "2.545.19"
$ws.Range("E16").Value = "  +1.80%  "

$ws.Range("E17").Value = "  -2.23%  "

$ws.Range("D18").Value = "42.771.46"
$ws.Range("E18").Value = "  -0.44%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.18"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.98%  "

$ws.Range("D20").Value = "0.0₃0968"
$ws.Range("E20").Value = "  -2.52%  "

$ws.Range("E21").Value = "  -0.85%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "71.28"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.98%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "251.94"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.96%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.98"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.26%  "

$ws.Range("E25").Value = "  -1.98%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "27.19"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.07%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.999"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.04%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.34"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +11.42%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "39.36"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.42%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "10.07"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.24%  "

$ws.Range("E31").Value = "  -3.83%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "156.94"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.08%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "19.76"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.97%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.32"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.43%  "

$ws.Range("E35").Value = "  -3.89%  "

$ws.Range("E36").Value = "  -2.25%  "

$ws.Range("E38").Value = "  -2.63%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "24.54"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -5.09%  "

$ws.Range("E40").Value = "  -0.30%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.11"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.19%  "

$ws.Range("E42").Value = "  -1.09%  "

$ws.Range("E43").Value = "  -1.45%  "

$ws.Range("D44").Value = "2.071.91"
$ws.Range("E44").Value = "  -0.51%  "

$ws.Range("E45").Value = "  +0.12%  "

$ws.Range("E46").Value = "  -1.88%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "86.88"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.42%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.82"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.98%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "74.42"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.60%  "

$ws.Range("D50").Value = "2.761.33"
$ws.Range("E50").Value = "  -1.45%  "

$ws.Range("E51").Value = "  -0.39%  "
